$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) updates
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 199
$wsExhibition.Range("F6").Value = 16
$wsExhibition.Range("F9").Value = 310
$wsExhibition.Range("F10").Value = 3207
$wsExhibition.Range("F11").Value = 31

# Sheet "全部类型" (All types) updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 199
$wsAll.Range("F7").Value = 16
$wsAll.Range("F10").Value = 310
$wsAll.Range("F11").Value = 3207
$wsAll.Range("F12").Value = 31
